$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The guild member row "Kyo" (row 35) was removed from the shard-payout
# sheet. Deleting the entire row shifts every row below it up by one,
# which Excel/the workbook automatically reflects in the dimension,
# merged cell ranges, conditional formatting ranges and any formulas
# that referenced rows below the deleted one (e.g. the H49/H51/H62/H67
# lookups in row 1 become H48/H50/H61/H66).
$ws.Rows.Item(35).Select()
$ws.Rows.Item(35).Delete()
